# Applies the crypto-price-tracker data refresh described by the commit diff.
# (GitHub Actions scheduled update of cryptos.xlsx - Sat Oct 26 15:24:09 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.895.84"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "2.456.56"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.86"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.84"
$ws.Range("E6").Value = "  -4.07%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("D9").Value = "2.451.49"
$ws.Range("E9").Value = "  -4.01%  "
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.87"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.326"
$ws.Range("E13").Value = "  -5.48%  "
$ws.Range("D14").Value = "2.900.11"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.08"
$ws.Range("E15").Value = "  -5.85%  "
$ws.Range("D16").Value = "66.807.06"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000167"
$ws.Range("E17").Value = "  -5.57%  "
$ws.Range("D18").Value = "2.397.60"
$ws.Range("E18").Value = "  -7.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  -8.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  -8.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.58"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.99"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.39"
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.18"
$ws.Range("E26").Value = "  -8.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.79"
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("E28").Value = "  -8.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -47.41%  "
$ws.Range("D30").Value = "2.576.51"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.03"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("D32").Value = "0.0₃0890"
$ws.Range("E32").Value = "  -7.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.56"
$ws.Range("E33").Value = "  -9.16%  "
$ws.Range("E34").Value = "  -6.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("E35").Value = "  -7.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.82"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  -12.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.63"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.12"
$ws.Range("E40").Value = "  -6.23%  "
$ws.Range("E41").Value = "  -9.00%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("E43").Value = "  -7.32%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.73"
$ws.Range("E44").Value = "  -8.10%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.322"
$ws.Range("E45").Value = "  -7.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.47"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.33"
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.42"
$ws.Range("E49").Value = "  -8.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.504"
$ws.Range("E50").Value = "  -9.35%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0248"
$ws.Range("E51").Value = "  -10.95%  "
